$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- original row 22
$ws.Range("D2").Value = 45163
$ws.Range("K2").Value = "Hayward"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 270
$ws.Range("N2").Value = 19000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 19500
$ws.Range("Q2").Value = "`$/bandeja 18 kilos"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1083
$ws.Range("T2").Value = 18

# Row 3 <- original row 6
$ws.Range("D3").Value = 45034
$ws.Range("K3").Value = "Hayward"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 25000
$ws.Range("O3").Value = 26000
$ws.Range("P3").Value = 25600
$ws.Range("Q3").Value = "`$/bandeja 18 kilos"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1422
$ws.Range("T3").Value = 18

# Row 4 <- original row 5
$ws.Range("D4").Value = 45069
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 370
$ws.Range("N4").Value = 19000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 19486
$ws.Range("Q4").Value = "`$/bandeja 18 kilos"
$ws.Range("R4").Value = "Región Metropolitana"
$ws.Range("S4").Value = 1083
$ws.Range("T4").Value = 18

# Row 5 <- original row 7
$ws.Range("D5").Value = 44291
$ws.Range("K5").Value = "Hayward"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 17000
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 17500
$ws.Range("Q5").Value = "`$/bandeja 18 kilos"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 972
$ws.Range("T5").Value = 18

# Row 6 <- original row 35
$ws.Range("D6").Value = 45148
$ws.Range("K6").Value = "Hayward"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 300
$ws.Range("N6").Value = 22000
$ws.Range("O6").Value = 23000
$ws.Range("P6").Value = 22500
$ws.Range("Q6").Value = "`$/bandeja 18 kilos"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1250
$ws.Range("T6").Value = 18

# Row 7 <- original row 33
$ws.Range("D7").Value = 45134
$ws.Range("K7").Value = "Hayward"
$ws.Range("L7").Value = "Especial"
$ws.Range("M7").Value = 350
$ws.Range("N7").Value = 21000
$ws.Range("O7").Value = 22000
$ws.Range("P7").Value = 21429
$ws.Range("Q7").Value = "`$/bandeja 18 kilos"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 1190
$ws.Range("T7").Value = 18

# Row 8 <- original row 30
$ws.Range("D8").Value = 45169
$ws.Range("K8").Value = "Hayward"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 270
$ws.Range("N8").Value = 27000
$ws.Range("O8").Value = 28000
$ws.Range("P8").Value = 27500
$ws.Range("Q8").Value = "`$/bandeja 18 kilos"
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 1528
$ws.Range("T8").Value = 18

# Row 9 <- original row 31
$ws.Range("D9").Value = 45169
$ws.Range("K9").Value = "Hayward"
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 150
$ws.Range("N9").Value = 26000
$ws.Range("O9").Value = 26000
$ws.Range("P9").Value = 26000
$ws.Range("Q9").Value = "`$/bandeja 18 kilos"
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 1444
$ws.Range("T9").Value = 18

# Row 10 <- original row 26
$ws.Range("D10").Value = 44487
$ws.Range("K10").Value = "Hayward"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 300
$ws.Range("N10").Value = 14000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 14500
$ws.Range("Q10").Value = "`$/bandeja 10 kilos"
$ws.Range("R10").Value = "Región de O'Higgins"
$ws.Range("S10").Value = 1450
$ws.Range("T10").Value = 10

# Row 11 <- original row 24
$ws.Range("D11").Value = 45127
$ws.Range("K11").Value = "Hayward"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = 21000
$ws.Range("O11").Value = 22000
$ws.Range("P11").Value = 21500
$ws.Range("Q11").Value = "`$/bandeja 18 kilos"
$ws.Range("R11").Value = "Región de O'Higgins"
$ws.Range("S11").Value = 1194
$ws.Range("T11").Value = 18

# Row 13 <- original row 25
$ws.Range("D13").Value = 45043
$ws.Range("K13").Value = "Hayward"
$ws.Range("L13").Value = "Segunda"
$ws.Range("M13").Value = 300
$ws.Range("N13").Value = 21000
$ws.Range("O13").Value = 22000
$ws.Range("P13").Value = 21500
$ws.Range("Q13").Value = "`$/bandeja 18 kilos"
$ws.Range("R13").Value = "Región de O'Higgins"
$ws.Range("S13").Value = 1194
$ws.Range("T13").Value = 18

# Row 14 <- original row 11
$ws.Range("D14").Value = 44819
$ws.Range("K14").Value = "Hayward"
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = 17000
$ws.Range("O14").Value = 18000
$ws.Range("P14").Value = 17500
$ws.Range("Q14").Value = "`$/bandeja 10 kilos"
$ws.Range("R14").Value = "Región de O'Higgins"
$ws.Range("S14").Value = 1750
$ws.Range("T14").Value = 10

# Row 15 <- original row 10
$ws.Range("D15").Value = 44418
$ws.Range("K15").Value = "Hayward"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 240
$ws.Range("N15").Value = 10000
$ws.Range("O15").Value = 11000
$ws.Range("P15").Value = 10500
$ws.Range("Q15").Value = "`$/bandeja 10 kilos"
$ws.Range("R15").Value = "Región de O'Higgins"
$ws.Range("S15").Value = 1050
$ws.Range("T15").Value = 10

# Row 16 <- original row 17
$ws.Range("D16").Value = 44263
$ws.Range("K16").Value = "Hayward"
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 250
$ws.Range("N16").Value = 21000
$ws.Range("O16").Value = 22000
$ws.Range("P16").Value = 21500
$ws.Range("Q16").Value = "`$/caja 18 kilos"
$ws.Range("R16").Value = "Región de O'Higgins"
$ws.Range("S16").Value = 1194
$ws.Range("T16").Value = 18

# Row 17 <- original row 19
$ws.Range("D17").Value = 44491
$ws.Range("K17").Value = "Hayward"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 300
$ws.Range("N17").Value = 14000
$ws.Range("O17").Value = 15000
$ws.Range("P17").Value = 14500
$ws.Range("Q17").Value = "`$/bandeja 10 kilos"
$ws.Range("R17").Value = "Región de O'Higgins"
$ws.Range("S17").Value = 1450
$ws.Range("T17").Value = 10

# Row 18 <- original row 15
$ws.Range("D18").Value = 44616
$ws.Range("K18").Value = "Hayward"
$ws.Range("L18").Value = "Segunda"
$ws.Range("M18").Value = 300
$ws.Range("N18").Value = 16000
$ws.Range("O18").Value = 17000
$ws.Range("P18").Value = 16500
$ws.Range("Q18").Value = "`$/caja 18 kilos granel"
$ws.Range("R18").Value = "Región de O'Higgins"
$ws.Range("S18").Value = 917
$ws.Range("T18").Value = 18

# Row 19 <- original row 18
$ws.Range("D19").Value = 44614
$ws.Range("K19").Value = "Hayward"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 250
$ws.Range("N19").Value = 20000
$ws.Range("O19").Value = 21000
$ws.Range("P19").Value = 20500
$ws.Range("Q19").Value = "`$/bandeja 18 kilos"
$ws.Range("R19").Value = "Región de O'Higgins"
$ws.Range("S19").Value = 1139
$ws.Range("T19").Value = 18

# Row 20 <- original row 8
$ws.Range("D20").Value = 44629
$ws.Range("K20").Value = "Hayward"
$ws.Range("L20").Value = "Segunda"
$ws.Range("M20").Value = 300
$ws.Range("N20").Value = 17000
$ws.Range("O20").Value = 18000
$ws.Range("P20").Value = 17500
$ws.Range("Q20").Value = "`$/bandeja 18 kilos"
$ws.Range("R20").Value = "Región de O'Higgins"
$ws.Range("S20").Value = 972
$ws.Range("T20").Value = 18

# Row 21 <- original row 27
$ws.Range("D21").Value = 44784
$ws.Range("K21").Value = "Hayward"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 300
$ws.Range("N21").Value = 19000
$ws.Range("O21").Value = 20000
$ws.Range("P21").Value = 19500
$ws.Range("Q21").Value = "`$/bandeja 18 kilos"
$ws.Range("R21").Value = "Región de O'Higgins"
$ws.Range("S21").Value = 1083
$ws.Range("T21").Value = 18

# Row 22 <- original row 28
$ws.Range("D22").Value = 44323
$ws.Range("K22").Value = "Hayward"
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 270
$ws.Range("N22").Value = 21000
$ws.Range("O22").Value = 22000
$ws.Range("P22").Value = 21500
$ws.Range("Q22").Value = "`$/bandeja 18 kilos"
$ws.Range("R22").Value = "Región de O'Higgins"
$ws.Range("S22").Value = 1194
$ws.Range("T22").Value = 18

# Row 24 <- original row 34
$ws.Range("D24").Value = 44489
$ws.Range("K24").Value = "Hayward"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 300
$ws.Range("N24").Value = 26000
$ws.Range("O24").Value = 27000
$ws.Range("P24").Value = 26500
$ws.Range("Q24").Value = "`$/bandeja 18 kilos"
$ws.Range("R24").Value = "Región de O'Higgins"
$ws.Range("S24").Value = 1472
$ws.Range("T24").Value = 18

# Row 25 <- original row 29
$ws.Range("D25").Value = 44673
$ws.Range("K25").Value = "Hayward"
$ws.Range("L25").Value = "Especial"
$ws.Range("M25").Value = 400
$ws.Range("N25").Value = 14000
$ws.Range("O25").Value = 15000
$ws.Range("P25").Value = 14500
$ws.Range("Q25").Value = "`$/bandeja 10 kilos"
$ws.Range("R25").Value = "Región de O'Higgins"
$ws.Range("S25").Value = 1450
$ws.Range("T25").Value = 10

# Row 26 <- original row 16
$ws.Range("D26").Value = 45204
$ws.Range("K26").Value = "Hayward"
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 300
$ws.Range("N26").Value = 23000
$ws.Range("O26").Value = 24000
$ws.Range("P26").Value = 23500
$ws.Range("Q26").Value = "`$/bandeja 10 kilos"
$ws.Range("R26").Value = "Región de O'Higgins"
$ws.Range("S26").Value = 2350
$ws.Range("T26").Value = 10

# Row 27 <- original row 20
$ws.Range("D27").Value = 44656
$ws.Range("K27").Value = "Hayward"
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 270
$ws.Range("N27").Value = 19000
$ws.Range("O27").Value = 20000
$ws.Range("P27").Value = 19500
$ws.Range("Q27").Value = "`$/bandeja 18 kilos"
$ws.Range("R27").Value = "Región de O'Higgins"
$ws.Range("S27").Value = 1083
$ws.Range("T27").Value = 18

# Row 28 <- original row 13
$ws.Range("D28").Value = 44991
$ws.Range("K28").Value = "Hayward"
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 250
$ws.Range("N28").Value = 24000
$ws.Range("O28").Value = 25000
$ws.Range("P28").Value = 24500
$ws.Range("Q28").Value = "`$/bandeja 18 kilos"
$ws.Range("R28").Value = "Región de O'Higgins"
$ws.Range("S28").Value = 1361
$ws.Range("T28").Value = 18

# Row 29 <- original row 14
$ws.Range("D29").Value = 44307
$ws.Range("K29").Value = "Hayward"
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 250
$ws.Range("N29").Value = 19000
$ws.Range("O29").Value = 20000
$ws.Range("P29").Value = 19500
$ws.Range("Q29").Value = "`$/bandeja 18 kilos"
$ws.Range("R29").Value = "Región de O'Higgins"
$ws.Range("S29").Value = 1083
$ws.Range("T29").Value = 18

# Row 30 <- original row 21
$ws.Range("D30").Value = 45223
$ws.Range("K30").Value = "Hayward"
$ws.Range("L30").Value = "Especial"
$ws.Range("M30").Value = 250
$ws.Range("N30").Value = 24000
$ws.Range("O30").Value = 25000
$ws.Range("P30").Value = 24600
$ws.Range("Q30").Value = "`$/bandeja 10 kilos"
$ws.Range("R30").Value = "Región de O'Higgins"
$ws.Range("S30").Value = 2460
$ws.Range("T30").Value = 10

# Row 31 <- original row 3
$ws.Range("D31").Value = 45086
$ws.Range("K31").Value = "Hayward"
$ws.Range("L31").Value = "Especial"
$ws.Range("M31").Value = 250
$ws.Range("N31").Value = 25000
$ws.Range("O31").Value = 26000
$ws.Range("P31").Value = 25500
$ws.Range("Q31").Value = "`$/bandeja 18 kilos"
$ws.Range("R31").Value = "Región de O'Higgins"
$ws.Range("S31").Value = 1417
$ws.Range("T31").Value = 18

# Row 32 <- original row 4
$ws.Range("D32").Value = 45086
$ws.Range("K32").Value = "Hayward"
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 250
$ws.Range("N32").Value = 20000
$ws.Range("O32").Value = 21000
$ws.Range("P32").Value = 20500
$ws.Range("Q32").Value = "`$/bandeja 18 kilos"
$ws.Range("R32").Value = "Región de O'Higgins"
$ws.Range("S32").Value = 1139
$ws.Range("T32").Value = 18

# Row 33 <- original row 9
$ws.Range("D33").Value = 44706
$ws.Range("K33").Value = "Hayward"
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 400
$ws.Range("N33").Value = 9000
$ws.Range("O33").Value = 10000
$ws.Range("P33").Value = 9500
$ws.Range("Q33").Value = "`$/bandeja 10 kilos"
$ws.Range("R33").Value = "Región de O'Higgins"
$ws.Range("S33").Value = 950
$ws.Range("T33").Value = 10

# Row 34 <- original row 32
$ws.Range("D34").Value = 44789
$ws.Range("K34").Value = "Hayward"
$ws.Range("L34").Value = "Segunda"
$ws.Range("M34").Value = 250
$ws.Range("N34").Value = 19000
$ws.Range("O34").Value = 20000
$ws.Range("P34").Value = 19500
$ws.Range("Q34").Value = "`$/bandeja 18 kilos"
$ws.Range("R34").Value = "Región de O'Higgins"
$ws.Range("S34").Value = 1083
$ws.Range("T34").Value = 18

# Row 35 <- original row 2
$ws.Range("D35").Value = 45002
$ws.Range("K35").Value = "Hayward"
$ws.Range("L35").Value = "Segunda"
$ws.Range("M35").Value = 300
$ws.Range("N35").Value = 24000
$ws.Range("O35").Value = 25000
$ws.Range("P35").Value = 24500
$ws.Range("Q35").Value = "`$/bandeja 18 kilos"
$ws.Range("R35").Value = "Región de O'Higgins"
$ws.Range("S35").Value = 1361
$ws.Range("T35").Value = 18
